$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 148 (weekly update), shifting rows 148:159 down to 149:160
$ws.Rows(148).Insert()

# Populate the newly inserted row 148 with a fresh record (same stall/product data as the
# previous row 148, but a new survey date - 2023-07-25)
$ws.Cells.Item(148, 1).Value2 = 8
$ws.Cells.Item(148, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(148, 3).Value2 = "Coquimbo"
$ws.Cells.Item(148, 4).Value2 = 45132
$ws.Cells.Item(148, 5).Value2 = 4
$ws.Cells.Item(148, 6).Value2 = 100114007
$ws.Cells.Item(148, 7).Value2 = "Jengibre"
$ws.Cells.Item(148, 8).Value2 = "Sin especificar"
$ws.Cells.Item(148, 9).Value2 = "Primera"
$ws.Cells.Item(148, 10).Value2 = 360
$ws.Cells.Item(148, 11).Value2 = 17000
$ws.Cells.Item(148, 12).Value2 = 18000
$ws.Cells.Item(148, 13).Value2 = 17500
$ws.Cells.Item(148, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(148, 15).Value2 = "Perú"
$ws.Cells.Item(148, 16).Value2 = 1346
$ws.Cells.Item(148, 17).Value2 = 13
$ws.Cells.Item(148, 18).Value2 = "Hortaliza"
